$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New koff reference values in column F for rows 2-5 (Stefan et al. 2012 / Pharris et al. 2020 update) ---
$ws.Range("F2").Value = 1.928
$ws.Range("F3").Value = 1.937
$ws.Range("F4").Value = 7.476
$ws.Range("F5").Value = 25.783

# --- New literature references for the CaMKII subunit active/inactive flicker row (row 6) ---
$ws.Range("D6").WrapText = $true
$ws.Range("D6").Value = "Stefan et al., 2012; Pharris et al., 2020"
$ws.Range("G6").Value = "Pharris et al., 2020"

# Row 6 grows to fit the wrapped, multi-line reference text
$ws.Rows.Item(6).RowHeight = 60

# --- New Kcat value for the CaMKII dephosphorylation by PP1 row (row 12) ---
$ws.Range("C2").Copy($ws.Range("C12"))
$ws.Range("C12").Value = 20

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Restore the active selection to match the saved view ---
$ws.Range("F11").Select()
